$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new time-log entry (row 16) ---------------------------------
# Date: 2017-12-17 (serial 43086)
$ws.Range("A16").Value = 43086
# User: Giovanni (existing shared string)
$ws.Range("B16").Value = "Giovanni"
# Activity: new note describing the work done
$ws.Range("C16").Value = "Continuato a scrivere codice. Classe oggetto riempita con tutti I parametri e le funzioni di servizio"
# Hours: 1 hour = 1/24 of a day
$ws.Range("D16").Value = 1/24

# The activity text wraps across several lines in a narrow column, so the
# row needs to grow to show it all (mirrors Excel's automatic row resize).
$ws.Rows.Item(16).RowHeight = 72.9

# --- Update the view: scroll down and select the next empty activity cell -
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
[void]$ws.Range("D17").Select()
